$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct tiny floating point drift on existing row 9, column A
$ws.Cells.Item(9, 1).Value = 45875.41708634259

# Append new row 10 with the new weather reading
$ws.Cells.Item(10, 1).Value = 45875.45855712878
$ws.Cells.Item(10, 2).Value = 2025
$ws.Cells.Item(10, 3).Value = 23
$ws.Cells.Item(10, 4).Value = 19.33
$ws.Cells.Item(10, 5).Value = 78.64
$ws.Cells.Item(10, 6).Value = 564.26
$ws.Cells.Item(10, 7).Value = 11.2
$ws.Cells.Item(10, 8).Value = "ESE"
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = "11:00:19"

# Match the date/time number format used by the date column in previous rows
$ws.Cells.Item(10, 1).NumberFormat = $ws.Cells.Item(9, 1).NumberFormat
